$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Row 7 currently holds "F Suite" in column A (TSID). Replace it with "Notifications".
$ws.Range("A7").Value = "Notifications"

# Update the selected cell to reflect the edited row.
$ws.Range("A7").Select()
